$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing "Fecha real de cierre" (actual closing date) values,
# matching the planned closing date already present in column D for rows 4 and 5.
$ws.Range("E4").Value = 42360
$ws.Range("E5").Value = 42360

# Update the active cell selection to C13 (cursor moved after edit)
$ws.Range("C13").Select()
